# FeatherFriend / BirdDB.xlsx update
# - changed database of genetics to match all body color cases
# - updated search cage to match the assignment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing BodyColor (column K) for cage/bird on row 6
$ws.Range("K6").Value = "Green"

# Add a new bird record on row 9 (genetics / color case addition)
$ws.Range("A9").Value = 207338351
$ws.Range("B9").Value = "American Gouldian"
$ws.Range("C9").Value = "North America"
$ws.Range("D9").Value = "453A"
$ws.Range("E9").Value = "Male"
$ws.Range("F9").Value = 223
$ws.Range("G9").Value = 435345
$ws.Range("H9").Value = "15/05/2023"
$ws.Range("I9").Value = "Red"
$ws.Range("J9").Value = "Purple"
$ws.Range("K9").Value = "White Pastel"

# Update the active selection to reflect the search/assignment cell used last
$ws.Range("L6").Select()
